# "added details on the example"
#
# The "simple" NPV worksheet had its discounted-cash-flow column (C4:C8)
# referencing an empty cell ($F$2) instead of the discount-rate cell
# ($B$2). Fix the formulas so the sheet actually computes discounted
# cash flows, copy the bottom-border formatting that row 8 already had
# in column B onto C8 to match, and leave the "simple" tab selected /
# scrolled to D20 as the last-saved view state (mirroring the source
# workbook, which also moved the active tab from "npv v2" to "simple").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("simple")

$ws.Range("C4").Formula = '=B4/((1+$B$2)^A4)'
$ws.Range("C5").Formula = '=B5/((1+$B$2)^A5)'
$ws.Range("C6").Formula = '=B6/((1+$B$2)^A6)'
$ws.Range("C7").Formula = '=B7/((1+$B$2)^A7)'
$ws.Range("C8").Formula = '=B8/((1+$B$2)^A8)'

# Row 8's cash-flow cell (B8) already carries a bottom-border style;
# copy that formatting across to the discounted-cash-flow cell (C8) so
# the two line up visually, without disturbing C8's formula/value.
$ws.Range("B8").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Make "simple" the active sheet / tab, with D20 as the selected cell.
$ws.Activate()
$ws.Range("D20").Select() | Out-Null
